$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4081.6365
$ws.Range("I51").Value = 3787.5
$ws.Range("J51").Value = 4866
$ws.Range("K51").Value = 3787.5
$ws.Range("L51").Value = 4866
$ws.Range("M51").Value = -3303.5
$ws.Range("N51").Value = -5834

$ws.Range("H98").Value = 1200
$ws.Range("I98").Value = 1200
$ws.Range("K98").Value = 1200
$ws.Range("M98").Value = 298

$ws.Range("H106").Value = 51900
$ws.Range("I106").Value = 51900
$ws.Range("K106").Value = 51900
$ws.Range("M106").Value = -51269

$ws.Range("H113").Value = 3866.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3866.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = ""
$ws.Range("M113").Value = 3866.5
$ws.Range("N113").Value = -10374.5

$ws.Range("H122").Value = 1200
$ws.Range("I122").Value = 1200
$ws.Range("K122").Value = 3600
$ws.Range("M122").Value = -1150

$ws.Range("H125").Value = 250015740
$ws.Range("J125").Value = 100018900
$ws.Range("L125").Value = 900170100
$ws.Range("N125").Value = -900175020

$ws.Range("H138").Value = 4877.8965
$ws.Range("I138").Value = 2684.3333
$ws.Range("J138").Value = 5865
$ws.Range("K138").Value = 8052.999899999999
$ws.Range("L138").Value = 17595
$ws.Range("M138").Value = -2912.999899999999
$ws.Range("N138").Value = -27875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3499
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""

$ws.Range("H45").Value = 7228.75
$ws.Range("I45").Value = 9133.333000000001
$ws.Range("J45").Value = 1515
$ws.Range("K45").Value = 9133.333000000001
$ws.Range("L45").Value = 1515
$ws.Range("M45").Value = -8756.333000000001
$ws.Range("N45").Value = -2269

$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 3000
$ws.Range("K61").Value = 3000
$ws.Range("M61").Value = -2788

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = ""
$ws.Range("N76").Value = 0

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = ""
$ws.Range("N79").Value = 0

$ws.Range("H116").Value = 3499
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").Value = ""

$ws.Range("H122").Value = 1818.6
$ws.Range("I122").Value = 1492.5294
$ws.Range("J122").Value = 3666.3333
$ws.Range("K122").Value = 4477.5882
$ws.Range("L122").Value = 10998.9999
$ws.Range("M122").Value = -2027.5882
$ws.Range("N122").Value = -15898.9999

$ws.Range("H128").Value = 100000
$ws.Range("J128").Value = 100000
$ws.Range("L128").Value = 100000
$ws.Range("N128").Value = -109960

$ws.Range("H132").Value = 2739.5
$ws.Range("I132").Value = 2883.4
$ws.Range("K132").Value = 8650.200000000001
$ws.Range("M132").Value = -6120.200000000001

$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -6450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3499
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""

$ws.Range("H86").Value = 1554.1818
$ws.Range("I86").Value = 1529.6
$ws.Range("K86").Value = 1529.6
$ws.Range("M86").Value = -406.5999999999999

$ws.Range("H89").Value = 1554.1818
$ws.Range("I89").Value = 1529.6
$ws.Range("K89").Value = 7648
$ws.Range("M89").Value = -2032

$ws.Range("H134").Value = 2560
$ws.Range("I134").Value = 2614.6
$ws.Range("K134").Value = 7843.799999999999
$ws.Range("M134").Value = -5308.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 70.72727
$ws.Range("I7").Value = 25.5
$ws.Range("J7").Value = 80.77778000000001
$ws.Range("K7").Value = 25.5
$ws.Range("L7").Value = 80.77778000000001
$ws.Range("M7").Value = 87.5
$ws.Range("N7").Value = -306.77778

$ws.Range("H16").Value = 6500
$ws.Range("I16").Value = 6500
$ws.Range("K16").Value = 6500
$ws.Range("M16").Value = -6213

$ws.Range("H22").Value = 497.5
$ws.Range("I22").Value = 497.5
$ws.Range("K22").Value = 497.5
$ws.Range("M22").Value = -147.5

$ws.Range("H58").Value = 2350.8572
$ws.Range("I58").Value = 2002
$ws.Range("K58").Value = 2002
$ws.Range("M58").Value = -1799

$ws.Range("H94").Value = 96197.164
$ws.Range("I94").Value = 161679.72
$ws.Range("J94").Value = 4521.6
$ws.Range("K94").Value = 161679.72
$ws.Range("L94").Value = 4521.6
$ws.Range("M94").Value = -161228.72
$ws.Range("N94").Value = -5423.6

$ws.Range("H113").Value = 6500
$ws.Range("I113").Value = 6500
$ws.Range("K113").Value = 6500
$ws.Range("M113").Value = -4330

$ws.Range("H132").Value = 4455.5557
$ws.Range("I132").Value = 4625
$ws.Range("K132").Value = 13875
$ws.Range("M132").Value = -11345

$ws.Range("H136").Value = 2350.8572
$ws.Range("I136").Value = 2002
$ws.Range("K136").Value = 6006
$ws.Range("M136").Value = -3456

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 499
$ws.Range("I51").Value = 499
$ws.Range("K51").Value = 1497
$ws.Range("M51").Value = -1037

$ws.Range("H81").Value = 2284.1428
$ws.Range("J81").Value = 2398
$ws.Range("L81").Value = 7194
$ws.Range("N81").Value = -9440

$ws.Range("H84").Value = 2284.1428
$ws.Range("J84").Value = 2398
$ws.Range("L84").Value = 21582
$ws.Range("N84").Value = -32814

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = ""

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = ""
$ws.Range("N21").Value = 0

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = ""
$ws.Range("N30").Value = 0

$ws.Range("H34").Value = 45666.332
$ws.Range("J34").Value = 49999.5
$ws.Range("L34").Value = 49999.5
$ws.Range("N34").Value = -50535.5

$ws.Range("H76").Value = 45666.332
$ws.Range("J76").Value = 49999.5
$ws.Range("L76").Value = 49999.5
$ws.Range("N76").Value = -50629.5

$ws.Range("H79").Value = 45666.332
$ws.Range("J79").Value = 49999.5
$ws.Range("L79").Value = 49999.5
$ws.Range("N79").Value = -52183.5

$ws.Range("H80").Value = 3866.6667
$ws.Range("I80").Value = 4200
$ws.Range("J80").Value = 3700
$ws.Range("K80").Value = 4200
$ws.Range("L80").Value = 3700
$ws.Range("M80").Value = -3202
$ws.Range("N80").Value = -5696

$ws.Range("H83").Value = 3866.6667
$ws.Range("I83").Value = 4200
$ws.Range("J83").Value = 3700
$ws.Range("K83").Value = 21000
$ws.Range("L83").Value = 18500
$ws.Range("M83").Value = -16008
$ws.Range("N83").Value = -28484

$ws.Range("H132").Value = 3924.2
$ws.Range("I132").Value = 3924.2
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11772.6
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = -9242.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1000000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = ""

$ws.Range("H46").Value = 3992.6667
$ws.Range("I46").Value = 3992.6667
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3992.6667
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -3804.6667

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = ""
$ws.Range("N51").Value = 0

$ws.Range("H55").Value = 1764.9474
$ws.Range("I55").Value = 1558.6364
$ws.Range("J55").Value = 2048.625
$ws.Range("K55").Value = 1558.6364
$ws.Range("L55").Value = 2048.625
$ws.Range("M55").Value = -1385.6364
$ws.Range("N55").Value = -2394.625

$ws.Range("H69").Value = 40000
$ws.Range("J69").Value = 40000
$ws.Range("L69").Value = 40000
$ws.Range("N69").Value = -41622

$ws.Range("H72").Value = 40000
$ws.Range("J72").Value = 40000
$ws.Range("L72").Value = 120000
$ws.Range("N72").Value = -128112

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19978.75
$ws.Range("I41").Value = 19978
$ws.Range("K41").Value = 19978
$ws.Range("M41").Value = -19588
